$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1932
$ws1.Range("F11").Value = 628
$ws1.Range("F14").Value = 698
$ws1.Range("F17").Value = 177

# Sheet: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1961

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1961
$ws4.Range("F16").Value = 1932
$ws4.Range("F26").Value = 628
$ws4.Range("F30").Value = 698
$ws4.Range("F35").Value = 177
